$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Increment the date serial values in column F (rows 2-7) by 1 day.
$ws.Range("F2").Value = 44482
$ws.Range("F3").Value = 44481
$ws.Range("F4").Value = 44480
$ws.Range("F5").Value = 44479
$ws.Range("F6").Value = 44478
$ws.Range("F7").Value = 44477
